# Auto update Excel log
# Appends newly-logged SeniorConnect sensor events to PIR, Humidity,
# Temperature, and Proximity sheets.

$wb = $excel.ActiveWorkbook

function Set-LogRow {
    param($ws, [int]$row, [string[]]$values)
    for ($c = 0; $c -lt $values.Count; $c++) {
        $cell = $ws.Cells.Item($row, $c + 1)
        $v = $values[$c]
        if ($v -match '^\d{4}-\d{2}-\d{2}$' -or $v -match '^\d+(\.\d+)?%$') {
            # Force text so Excel does not reinterpret dates/percentages as numbers
            $cell.NumberFormat = "@"
        }
        $cell.Value = $v
    }
}

$ws = $wb.Worksheets.Item("PIR")
Set-LogRow $ws 48 @('2026-02-01', '18:29:15', '18:00', 'Bathroom', 'No Motion', 'Inactive')
Set-LogRow $ws 49 @('2026-02-01', '18:29:15', '18:00', 'Bathroom', 'Motion Detected', 'Active')
Set-LogRow $ws 50 @('2026-02-01', '18:29:23', '18:00', 'Bathroom', 'No Motion', 'Inactive')
Set-LogRow $ws 51 @('2026-02-01', '18:29:25', '18:00', 'Bathroom', 'Motion Detected', 'Active')
Set-LogRow $ws 52 @('2026-02-01', '18:29:32', '18:00', 'Bathroom', 'No Motion', 'Inactive')
Set-LogRow $ws 53 @('2026-02-01', '18:29:33', '18:00', 'Bathroom', 'Motion Detected', 'Active')
Set-LogRow $ws 54 @('2026-02-01', '18:29:41', '18:00', 'Bathroom', 'No Motion', 'Inactive')

$ws = $wb.Worksheets.Item("Humidity")
Set-LogRow $ws 104 @('2026-02-01', '18:28:49', '18:00', 'Bathroom', '78.1%', 'Active')
Set-LogRow $ws 105 @('2026-02-01', '18:28:59', '18:00', 'Bathroom', '78.2%', 'Active')
Set-LogRow $ws 106 @('2026-02-01', '18:29:09', '18:00', 'Bathroom', '80.6%', 'Active')
Set-LogRow $ws 107 @('2026-02-01', '18:29:14', '18:00', 'Bathroom', '81.4%', 'Active')
Set-LogRow $ws 108 @('2026-02-01', '18:29:19', '18:00', 'Bathroom', '80.5%', 'Active')
Set-LogRow $ws 109 @('2026-02-01', '18:29:24', '18:00', 'Bathroom', '79.8%', 'Active')
Set-LogRow $ws 110 @('2026-02-01', '18:29:29', '18:00', 'Bathroom', '78.6%', 'Active')
Set-LogRow $ws 111 @('2026-02-01', '18:29:34', '18:00', 'Bathroom', '79.7%', 'Active')
Set-LogRow $ws 112 @('2026-02-01', '18:29:39', '18:00', 'Bathroom', '78.6%', 'Active')

$ws = $wb.Worksheets.Item("Temperature")
Set-LogRow $ws 104 @('2026-02-01', '18:28:49', '18:00', 'Bathroom', '29.6C', 'Active')
Set-LogRow $ws 105 @('2026-02-01', '18:28:59', '18:00', 'Bathroom', '29.6C', 'Active')
Set-LogRow $ws 106 @('2026-02-01', '18:29:09', '18:00', 'Bathroom', '29.5C', 'Active')
Set-LogRow $ws 107 @('2026-02-01', '18:29:14', '18:00', 'Bathroom', '29.6C', 'Active')
Set-LogRow $ws 108 @('2026-02-01', '18:29:19', '18:00', 'Bathroom', '29.6C', 'Active')
Set-LogRow $ws 109 @('2026-02-01', '18:29:24', '18:00', 'Bathroom', '29.6C', 'Active')
Set-LogRow $ws 110 @('2026-02-01', '18:29:30', '18:00', 'Bathroom', '29.6C', 'Active')
Set-LogRow $ws 111 @('2026-02-01', '18:29:35', '18:00', 'Bathroom', '29.6C', 'Active')
Set-LogRow $ws 112 @('2026-02-01', '18:29:40', '18:00', 'Bathroom', '29.6C', 'Active')

$ws = $wb.Worksheets.Item("Proximity")
Set-LogRow $ws 43 @('2026-02-01', '18:28:51', '18:00', 'Bathroom Door', 'EXIT', 'User EXITED Bathroom')

